# "!ref make me angy" - blow out the marks column with -8, except two rows
# that get their own special (still wrong) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most of the "marks" column (C2:C58 and C61:C67) gets slammed to -8.
$ws.Range("C2:C58").Value = -8
$ws.Range("C61:C67").Value = -8

# Two rows escape the -8 treatment and get different bogus values instead.
$ws.Range("C59").Value = 32
$ws.Range("C60").Value = 43
